$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row 1 cell text: "physical quantity" columns renamed to "measurand"
# (order matters for shared-string table append order, matching the target)
$ws.Range("I1").Value = "Measurand Level I"
$ws.Range("J1").Value = "Measurand Level II"
$ws.Range("F1").Value = "Measurand"

# Update sheet view: scroll/selection moved from H1/I2:I3 to F1
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("F1").Select()
